$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.86"
$ws.Range("E2").Value = "'0.75%"
$ws.Range("D3").Value = "'26.19"
$ws.Range("E3").Value = "'4.31%"
$ws.Range("D4").Value = "'5.080"
$ws.Range("E4").Value = "'1.58%"
$ws.Range("D6").Value = "'6.483"
$ws.Range("E6").Value = "'-1.41%"
$ws.Range("D7").Value = "'0.8128"
$ws.Range("E7").Value = "'0.20%"
$ws.Range("D8").Value = "'0.8431"
$ws.Range("E8").Value = "'-0.08%"
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").Value = "'0.07002"
$ws.Range("E9").Value = "'1.05%"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Value = "'0.02825"
$ws.Range("E10").Value = "'-0.49%"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Value = "'0.09398"
$ws.Range("E11").Value = "'-0.05%"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Value = "'0.001508"
$ws.Range("E12").Value = "'-0.52%"
$ws.Range("B13").Value = "TigerCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D13").Value = "'0.006153"
$ws.Range("E13").Value = "'1.10%"
$ws.Range("B14").Value = "LEO"
$ws.Range("C14").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D14").Value = "'3.607"
$ws.Range("E14").Value = "'3.09%"
$ws.Range("B15").Value = "GateToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D15").Value = "'3.013"
$ws.Range("E15").Value = "'0.17%"
$ws.Range("B16").Value = "BTSEToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D16").Value = "'2.055"
$ws.Range("E16").Value = "'-1.73%"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005948"
$ws.Range("E17").Value = "'-0.17%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3126"
$ws.Range("E18").Value = "'-1.27%"
$ws.Range("B19").Value = "WazirX"
$ws.Range("C19").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D19").Value = "'0.1338"
$ws.Range("E19").Value = "'0.02%"
$ws.Range("D20").Value = "'0.03204"
$ws.Range("E20").Value = "'-2.20%"
$ws.Range("D21").Value = "'0.1297"
$ws.Range("E21").Value = "'-1.67%"
$ws.Range("D22").Value = "'3.756"
$ws.Range("E22").Value = "'0.55%"
$ws.Range("D23").Value = "'0.04646"
$ws.Range("E23").Value = "'-0.49%"
$ws.Range("E24").Value = "'-1.46%"
$ws.Range("E25").Value = "'0.34%"
$ws.Range("D26").Value = "'0.004585"
$ws.Range("E26").Value = "'1.37%"
$ws.Range("E27").Value = "'-1.00%"
$ws.Range("D28").Value = "'0.0001937"
$ws.Range("E28").Value = "'-0.11%"
$ws.Range("D40").Value = "'0.03665"
$ws.Range("E40").Value = "'0.09%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006168"
$ws.Range("E41").Value = "'82.03%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1060"
$ws.Range("E42").Value = "'-21.87%"
$ws.Range("D43").Value = "'0.002499"
$ws.Range("E43").Value = "'-8.18%"
$ws.Range("D44").Value = "'0.008268"
$ws.Range("E44").Value = "'2.33%"
$ws.Range("D45").Value = "'0.00005387"
$ws.Range("E45").Value = "'1.76%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("E47").Value = "'-42.11%"
$ws.Range("D48").Value = "'0.002593"
$ws.Range("E48").Value = "'27.23%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.00%"
